$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 352, shifting existing rows 352:428 down to 353:429
$ws.Rows.Item(352).Insert()

# Populate the newly inserted row 352 with the new weekly price entry
$ws.Cells.Item(352, 1).Value = 4
$ws.Cells.Item(352, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(352, 3).Value = "Los Lagos"
$ws.Cells.Item(352, 4).Value = 45015
$ws.Cells.Item(352, 5).Value = 10
$ws.Cells.Item(352, 6).Value = 100112037
$ws.Cells.Item(352, 7).Value = "Cebollín"
$ws.Cells.Item(352, 8).Value = "Sin especificar"
$ws.Cells.Item(352, 9).Value = "Primera"
$ws.Cells.Item(352, 10).Value = 70
$ws.Cells.Item(352, 11).Value = 6500
$ws.Cells.Item(352, 12).Value = 7000
$ws.Cells.Item(352, 13).Value = 6750
$ws.Cells.Item(352, 14).Value = "`$/paquete 36 unidades"
$ws.Cells.Item(352, 15).Value = "Región Metropolitana"
$ws.Cells.Item(352, 16).Value = 188
$ws.Cells.Item(352, 17).Value = 36
$ws.Cells.Item(352, 18).Value = "Hortaliza"
